$d = $word.ActiveDocument
$t = $d.Tables.Item(2)

function Split-CategoryRun($RowNum, $FirstText, $FirstPreserve, $SecondText) {
    $cell = $t.Rows.Item($RowNum).Cells.Item(2)
    $p = $cell.Range.Paragraphs.Item(1)
    $pRange = $p.Range
    # Exclude the trailing paragraph mark from the replaced range.
    $full = $d.Range($pRange.Start, $pRange.End - 1)

    if ($FirstPreserve) {
        $firstRun = '<w:r><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="es-PA"/></w:rPr><w:t xml:space="preserve">' + $FirstText + '</w:t></w:r>'
    } else {
        $firstRun = '<w:r><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="es-PA"/></w:rPr><w:t>' + $FirstText + '</w:t></w:r>'
    }
    $secondRun = '<w:r><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="es-PA"/></w:rPr><w:t>' + $SecondText + '</w:t></w:r>'

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p><w:pPr><w:spacing w:before="120" w:after="120"/><w:rPr><w:b/><w:sz w:val="20"/><w:lang w:val="es-PA"/></w:rPr></w:pPr>' +
        $firstRun + $secondRun +
        '</w:p>' +
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $full.InsertXML($xml)
}

# R5.2 -> "Requisito de usabilidad"
Split-CategoryRun 8 "Requisito " $true "de usabilidad"
# R5.5 -> "Requisito de fiabilidad"
Split-CategoryRun 11 "Requisito " $true "de fiabilidad"
# R5.16 -> "Requisito de usabilidad"
Split-CategoryRun 22 "Requisito " $true "de usabilidad"
# R5.25 -> "Requisito de fiabilidad"
Split-CategoryRun 31 "Requisito " $true "de fiabilidad"
# R5.26 -> "Requisito de usabilidad"
Split-CategoryRun 32 "Requisito " $true "de usabilidad"
# R5.32 -> "Requisito de usabilidad"
Split-CategoryRun 38 "Requisito " $true "de usabilidad"
# R6 -> "Requisito del proyecto"
Split-CategoryRun 39 "Requisito de" $false "l proyecto"
